# Update registration summary numbers on the "Inscricoes" sheet.
# Diff summary:
#   Row 12: E12 28->29, F12 11->12, H12 13->14
#   Row 13: E13 5->6
#   Row 15: E15 93->95
#   Row 17: E17 21->22
#   Row 18: E18 92->94

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E12").Value = 29
$ws.Range("F12").Value = 12
$ws.Range("H12").Value = 14

$ws.Range("E13").Value = 6

$ws.Range("E15").Value = 95

$ws.Range("E17").Value = 22

$ws.Range("E18").Value = 94
